$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 44.855544
$ws.Range("H2").Value = 134.566632
$ws.Range("I2").Value = 0.9269627513664965
$ws.Range("J2").Value = 0.9269627513664968
$ws.Range("M2").Value = 0.3331563333333333
$ws.Range("N2").Value = 0.9994690000000001
$ws.Range("O2").Value = 0.02641712724508605
$ws.Range("P2").Value = 0.02641712724508605
$ws.Range("Q2").Value = 14.943908568712
$ws.Range("R2").Value = 134.495177118408
$ws.Range("S2").Value = 0.0244876929543038
$ws.Range("T2").Value = 0.0244876929543038

# Row 3
$ws.Range("G3").Value = 44.855544
$ws.Range("H3").Value = 134.566632
$ws.Range("I3").Value = 0.9269627513664965
$ws.Range("J3").Value = 0.9269627513664968
$ws.Range("O3").Value = 0.5158522877584945
$ws.Range("P3").Value = 0.5158522877584945
$ws.Range("Q3").Value = 291.812555987056
$ws.Range("R3").Value = 2626.313003883504
$ws.Range("S3").Value = 0.4781758559593157
$ws.Range("T3").Value = 0.4781758559593158

# Row 4
$ws.Range("G4").Value = 44.855544
$ws.Range("H4").Value = 134.566632
$ws.Range("I4").Value = 0.9269627513664965
$ws.Range("J4").Value = 0.9269627513664968
$ws.Range("M4").Value = 4.637142333333333
$ws.Range("N4").Value = 13.911427
$ws.Range("O4").Value = 0.3676951833620909
$ws.Range("P4").Value = 0.3676951833620909
$ws.Range("Q4").Value = 208.001541967096
$ws.Range("R4").Value = 1872.013877703864
$ws.Range("S4").Value = 0.3408397388335322
$ws.Range("T4").Value = 0.3408397388335323

# Row 5
$ws.Range("G5").Value = 44.855544
$ws.Range("H5").Value = 134.566632
$ws.Range("I5").Value = 0.9269627513664965
$ws.Range("J5").Value = 0.9269627513664968
$ws.Range("M5").Value = 1.135470333333333
$ws.Range("N5").Value = 3.406411
$ws.Range("O5").Value = 0.09003540163432863
$ws.Range("P5").Value = 0.09003540163432863
$ws.Range("Q5").Value = 50.93213949752801
$ws.Range("R5").Value = 458.389255477752
$ws.Range("S5").Value = 0.08345946361934482
$ws.Range("T5").Value = 0.08345946361934485

# Row 6
$ws.Range("I6").Value = 0.03026428998407557
$ws.Range("J6").Value = 0.03026428998407558
$ws.Range("M6").Value = 0.3331563333333333
$ws.Range("N6").Value = 0.9994690000000001
$ws.Range("O6").Value = 0.02641712724508605
$ws.Range("P6").Value = 0.02641712724508605
$ws.Range("Q6").Value = 0.4879017865090001
$ws.Range("R6").Value = 4.391116078581001
$ws.Range("S6").Value = 0.0007994955994915075
$ws.Range("T6").Value = 0.0007994955994915078

# Row 7
$ws.Range("I7").Value = 0.03026428998407557
$ws.Range("J7").Value = 0.03026428998407558
$ws.Range("O7").Value = 0.5158522877584945
$ws.Range("P7").Value = 0.5158522877584945
$ws.Range("S7").Value = 0.01561190322567187
$ws.Range("T7").Value = 0.01561190322567188

# Row 8
$ws.Range("I8").Value = 0.03026428998407557
$ws.Range("J8").Value = 0.03026428998407558
$ws.Range("M8").Value = 4.637142333333333
$ws.Range("N8").Value = 13.911427
$ws.Range("O8").Value = 0.3676951833620909
$ws.Range("P8").Value = 0.3676951833620909
$ws.Range("Q8").Value = 6.791016115747001
$ws.Range("R8").Value = 61.119145041723
$ws.Range("S8").Value = 0.01112803365501816
$ws.Range("T8").Value = 0.01112803365501816

# Row 9
$ws.Range("I9").Value = 0.03026428998407557
$ws.Range("J9").Value = 0.03026428998407558
$ws.Range("M9").Value = 1.135470333333333
$ws.Range("N9").Value = 3.406411
$ws.Range("O9").Value = 0.09003540163432863
$ws.Range("P9").Value = 0.09003540163432863
$ws.Range("Q9").Value = 1.662877000171
$ws.Range("R9").Value = 14.965893001539
$ws.Range("S9").Value = 0.002724857503894034
$ws.Range("T9").Value = 0.002724857503894034

# Row 10
$ws.Range("G10").Value = 1.967437666666666
$ws.Range("H10").Value = 5.902312999999999
$ws.Range("I10").Value = 0.04065810533109158
$ws.Range("J10").Value = 0.0406581053310916
$ws.Range("M10").Value = 0.3331563333333333
$ws.Range("N10").Value = 0.9994690000000001
$ws.Range("O10").Value = 0.02641712724508605
$ws.Range("P10").Value = 0.02641712724508605
$ws.Range("Q10").Value = 0.6554643190885555
$ws.Range("R10").Value = 5.899178871797
$ws.Range("S10").Value = 0.001074070342075558
$ws.Range("T10").Value = 0.001074070342075558

# Row 11
$ws.Range("G11").Value = 1.967437666666666
$ws.Range("H11").Value = 5.902312999999999
$ws.Range("I11").Value = 0.04065810533109158
$ws.Range("J11").Value = 0.0406581053310916
$ws.Range("O11").Value = 0.5158522877584945
$ws.Range("P11").Value = 0.5158522877584945
$ws.Range("Q11").Value = 12.79937691214289
$ws.Range("R11").Value = 115.194392209286
$ws.Range("S11").Value = 0.02097357665096943
$ws.Range("T11").Value = 0.02097357665096944

# Row 12
$ws.Range("G12").Value = 1.967437666666666
$ws.Range("H12").Value = 5.902312999999999
$ws.Range("I12").Value = 0.04065810533109158
$ws.Range("J12").Value = 0.0406581053310916
$ws.Range("M12").Value = 4.637142333333333
$ws.Range("N12").Value = 13.911427
$ws.Range("O12").Value = 0.3676951833620909
$ws.Range("P12").Value = 0.3676951833620909
$ws.Range("Q12").Value = 9.123288492294554
$ws.Range("R12").Value = 82.109596430651
$ws.Range("S12").Value = 0.01494978949487093
$ws.Range("T12").Value = 0.01494978949487093

# Row 13
$ws.Range("G13").Value = 1.967437666666666
$ws.Range("H13").Value = 5.902312999999999
$ws.Range("I13").Value = 0.04065810533109158
$ws.Range("J13").Value = 0.0406581053310916
$ws.Range("M13").Value = 1.135470333333333
$ws.Range("N13").Value = 3.406411
$ws.Range("O13").Value = 0.09003540163432863
$ws.Range("P13").Value = 0.09003540163432863
$ws.Range("Q13").Value = 2.233967103182555
$ws.Range("R13").Value = 20.105703928643
$ws.Range("S13").Value = 0.003660668843175669
$ws.Range("T13").Value = 0.00366066884317567

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1023373333333333
$ws.Range("H14").Value = 0.307012
$ws.Range("I14").Value = 0.002114853318336234
$ws.Range("J14").Value = 0.002114853318336234
$ws.Range("M14").Value = 0.3331563333333333
$ws.Range("N14").Value = 0.9994690000000001
$ws.Range("O14").Value = 0.02641712724508605
$ws.Range("P14").Value = 0.02641712724508605
$ws.Range("Q14").Value = 0.03409433073644445
$ws.Range("R14").Value = 0.306848976628
$ws.Range("S14").Value = 0.00005586834921518076
$ws.Range("T14").Value = 0.00005586834921518077

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1023373333333333
$ws.Range("H15").Value = 0.307012
$ws.Range("I15").Value = 0.002114853318336234
$ws.Range("J15").Value = 0.002114853318336234
$ws.Range("O15").Value = 0.5158522877584945
$ws.Range("P15").Value = 0.5158522877584945
$ws.Range("Q15").Value = 0.6657665062071111
$ws.Range("R15").Value = 5.991898555863999
$ws.Range("S15").Value = 0.00109095192253739
$ws.Range("T15").Value = 0.00109095192253739

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1023373333333333
$ws.Range("H16").Value = 0.307012
$ws.Range("I16").Value = 0.002114853318336234
$ws.Range("J16").Value = 0.002114853318336234
$ws.Range("M16").Value = 4.637142333333333
$ws.Range("N16").Value = 13.911427
$ws.Range("O16").Value = 0.3676951833620909
$ws.Range("P16").Value = 0.3676951833620909
$ws.Range("Q16").Value = 0.4745527806804444
$ws.Range("R16").Value = 4.270975026124
$ws.Range("S16").Value = 0.000777621378669568
$ws.Range("T16").Value = 0.0007776213786695681

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1023373333333333
$ws.Range("H17").Value = 0.307012
$ws.Range("I17").Value = 0.002114853318336234
$ws.Range("J17").Value = 0.002114853318336234
$ws.Range("M17").Value = 1.135470333333333
$ws.Range("N17").Value = 3.406411
$ws.Range("O17").Value = 0.09003540163432863
$ws.Range("P17").Value = 0.09003540163432863
$ws.Range("Q17").Value = 0.1162010059924445
$ws.Range("R17").Value = 1.045809053932
$ws.Range("S17").Value = 0.0001904116679140955
$ws.Range("T17").Value = 0.0001904116679140955
